# Maçı İzle (Watch Match) feature: add a "Link" column (H) with YouTube links
# for played matches, and fix the "Takım1"/"Takım2" header typo (Turkish "ı").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Maçlar")

# --- Header fixes -----------------------------------------------------
$ws.Range("D1").Value = "Takım1"
$ws.Range("E1").Value = "Takım2"
$ws.Range("H1").Value = "Link"

# --- New "Link" column width (matches the other custom column widths) --
$ws.Columns.Item(8).ColumnWidth = 40

# --- Link values for matches that have already been played -------------
$ws.Range("H2").Value  = "https://youtu.be/uAyL-1RXy30"
$ws.Range("H3").Value  = "https://youtu.be/HZrAxbLTD4E"
$ws.Range("H4").Value  = "https://youtu.be/xkq854V1MVs"
$ws.Range("H5").Value  = "https://youtu.be/nqf2RPCnLGw"
$ws.Range("H6").Value  = "https://youtu.be/OUcjYh9v4l0"
$ws.Range("H7").Value  = "https://youtu.be/0CQJbBd-DCo"
$ws.Range("H8").Value  = "https://youtu.be/Mfes4jxAMQ8"
$ws.Range("H9").Value  = "https://youtu.be/A4Zsh-nURAQ"
$ws.Range("H10").Value = "https://youtu.be/WNBXhqiaTQQ"
$ws.Range("H11").Value = "https://youtu.be/4udqsX0EYMA"

# Matches 12-19 have not been played yet / have no video link yet, so the
# Link cell is left blank (no button is rendered for them).

# --- Selection as left by the author after the edit ---------------------
$ws.Range("H2").Select()
